# Update "想去人数" (interest count) values that changed between scrapes.
# Sheet "展览": F3 541->542, F9 392->395, F10 3429->3438
# Sheet "全部类型": F4 541->542, F10 392->395, F11 3429->3438

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 542
$wsExhibit.Range("F9").Value = 395
$wsExhibit.Range("F10").Value = 3438

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 542
$wsAll.Range("F10").Value = 395
$wsAll.Range("F11").Value = 3438
